$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Confirmation Events")

# New "instructions" header in column C (copy the header formatting first)
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "instructions"

# Row 2 now holds "Attend Retreat" (previously on row 3) plus its instructions
$ws.Range("A2").Value = "Attend Retreat"
$ws.Range("B2").Value = 42493
$ws.Range("A1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "<h1>a heading</h1>`n<ul>`n<li>step 1</li>`n<li>step 2</li>`n</ul>`n<p> </p>`n<p> </p>"
$ws.Range("C2").WrapText = $true

# Row 3 now holds "Parent Information Meeting" (previously on row 2) plus its instructions
$ws.Range("A3").Value = "Parent Information Meeting"
$ws.Range("B3").Value = 42524
$ws.Range("A1").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "<p><em><strong>simple text</strong></em></p>"

$excel.CutCopyMode = $false

# Widen the instructions column and grow row 2 so the wrapped html is visible
$ws.Columns.Item(3).ColumnWidth = 77.45746071428571
$ws.Rows.Item(2).RowHeight = 86.6
$ws.Rows.Item(3).RowHeight = 14.6
